$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Title (A1): rename application name
# ---------------------------------------------------------------
$ws.Range("A1").Value = "Applicazione: Progetto-Angelo - Confronto Robustezza Locatori"

# ---------------------------------------------------------------
# Rows 4-9 ("LLM" block): fill in real numbers instead of the
# placeholder "N\D" text, apply the new center+vcenter style to
# columns D:F, and compute the failure-rate formula in column G
# (shared across G4:G10, matching the existing G11:G16 pattern).
# ---------------------------------------------------------------
$ws.Range("C4").Value = 38
$ws.Range("D4").Value = 32
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 5

$ws.Range("C5").Value = 38
$ws.Range("D5").Value = 20
$ws.Range("E5").Value = 13
$ws.Range("F5").Value = 5

$ws.Range("C6").Value = 38
$ws.Range("D6").Value = 27
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 5

$ws.Range("C7").Value = 38
$ws.Range("D7").Value = 31
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 5

$ws.Range("C8").Value = 38
$ws.Range("D8").Value = 27
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 5

$ws.Range("C9").Value = 38
$ws.Range("D9").Value = 31
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 5

# Style: columns D:F on rows 4-9 become centered + vertically centered
$ws.Range("D4:F9").VerticalAlignment = -4108

# Column G4:G10 gets the failure-rate formula (shared formula) with
# the same number format (0.00) + centered style already used by
# G11:G16 (style index 2).
$ws.Range("G4:G10").Formula = "=(E4/C4)*100"
$ws.Range("G4:G10").NumberFormat = "0.00"

# ---------------------------------------------------------------
# Row 19: add headers for the three new "generazioni" columns
# ---------------------------------------------------------------
$ws.Range("E19").Value = "Generazioni mancanti"
$ws.Range("F19").Value = "Generazioni non necessarie"
$ws.Range("G19").Value = "Generazioni necessarie ma errate"

# ---------------------------------------------------------------
# Row 20 ("LLM" totals): turn the placeholder text into real
# SUM formulas and fill in the new "generazioni" counts.
# ---------------------------------------------------------------
$ws.Range("B20").Formula = "=SUM(E4:E9)"
$ws.Range("C20").Formula = "=SUM(F4:F9)"
$ws.Range("D20").Formula = "=SUM(B20,C20)"
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 0

# ---------------------------------------------------------------
# Row 21 ("Analitica" totals): fill in the new "generazioni" counts
# ---------------------------------------------------------------
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0

# ---------------------------------------------------------------
# Column G width (approximate best-fit width for the new, longer
# header text in G19)
# ---------------------------------------------------------------
$ws.Columns("G").ColumnWidth = 27.35

# ---------------------------------------------------------------
# Selection moves to E25
# ---------------------------------------------------------------
$ws.Range("E25").Select() | Out-Null
